$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows to append (update from MV data feed).
$newRows = @(
    @{ Row = 282; Date = "08-10-2021"; B = 855; C = 4869; D = 73 },
    @{ Row = 283; Date = "09-10-2021"; B = 855; C = 4869; D = 73 },
    @{ Row = 284; Date = "10-10-2021"; B = 855; C = 4869; D = 73 },
    @{ Row = 285; Date = "11-10-2021"; B = 855; C = 4869; D = 73 },
    @{ Row = 286; Date = "12-10-2021"; B = 849; C = 4830; D = 72 },
    @{ Row = 287; Date = "13-10-2021"; B = 841; C = 4787; D = 71 }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row

    # Write the date as a formula that evaluates to the literal text, then
    # convert it to a plain value in place. This avoids Excel's "looks like
    # a date" auto-conversion (which would store a date serial + apply a
    # date number format) while also avoiding the quotePrefix style that a
    # leading-apostrophe text entry would leave behind. The net result is a
    # plain shared-string cell with no extra style, matching a date column
    # whose other rows are stored the same way.
    $dateCell = $ws.Range("A$rowNum")
    $dateCell.Formula = '="' + $r.Date + '"'
    $dateCell.Copy()
    $dateCell.PasteSpecial(-4163)
    $excel.CutCopyMode = $false

    $ws.Range("B$rowNum").Value = $r.B
    $ws.Range("C$rowNum").Value = $r.C
    $ws.Range("D$rowNum").Value = $r.D
}
